$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Addr="D2"; Val="23.457.09"},
    @{Addr="E2"; Val="  -1.33%  "},
    @{Addr="D3"; Val="1.637.38"},
    @{Addr="E3"; Val="  -1.19%  "},
    @{Addr="D4"; Val="'1.001"},
    @{Addr="E4"; Val="  +0.33%  "},
    @{Addr="D5"; Val="'1.000"},
    @{Addr="E5"; Val="  +0.21%  "},
    @{Addr="D6"; Val="'303.95"},
    @{Addr="E6"; Val="  -1.01%  "},
    @{Addr="D7"; Val="'0.3782"},
    @{Addr="E7"; Val="  -0.03%  "},
    @{Addr="D8"; Val="'51.89"},
    @{Addr="E8"; Val="  -1.82%  "},
    @{Addr="D9"; Val="'0.3618"},
    @{Addr="E9"; Val="  -1.93%  "},
    @{Addr="D10"; Val="'0.08186"},
    @{Addr="E10"; Val="  -0.03%  "},
    @{Addr="D11"; Val="'1.233"},
    @{Addr="E11"; Val="  -3.75%  "},
    @{Addr="D12"; Val="'1.001"},
    @{Addr="E12"; Val="  +0.36%  "},
    @{Addr="D13"; Val="'22.51"},
    @{Addr="E13"; Val="  -3.45%  "},
    @{Addr="D14"; Val="'6.455"},
    @{Addr="E14"; Val="  -4.40%  "},
    @{Addr="D15"; Val="'7.374"},
    @{Addr="E15"; Val="  -0.81%  "},
    @{Addr="D16"; Val="'0.00001241"},
    @{Addr="E16"; Val="  -2.94%  "},
    @{Addr="D17"; Val="1.632.06"},
    @{Addr="E17"; Val="  -1.69%  "},
    @{Addr="D18"; Val="'95.22"},
    @{Addr="E18"; Val="  -0.33%  "},
    @{Addr="D19"; Val="'0.06935"},
    @{Addr="E19"; Val="  +0.22%  "},
    @{Addr="D20"; Val="'6.591"},
    @{Addr="E20"; Val="  -0.68%  "},
    @{Addr="D21"; Val="'17.46"},
    @{Addr="E21"; Val="  -5.94%  "},
    @{Addr="D22"; Val="'1.001"},
    @{Addr="E22"; Val="  +0.28%  "},
    @{Addr="D23"; Val="'12.54"},
    @{Addr="E23"; Val="  -3.93%  "},
    @{Addr="D24"; Val="23.457.18"},
    @{Addr="E24"; Val="  -1.35%  "},
    @{Addr="D25"; Val="'2.522"},
    @{Addr="E25"; Val="  +3.89%  "},
    @{Addr="D26"; Val="'3.056"},
    @{Addr="E26"; Val="  -5.94%  "},
    @{Addr="D27"; Val="'21.15"},
    @{Addr="E27"; Val="  -1.71%  "},
    @{Addr="D28"; Val="'151.66"},
    @{Addr="E28"; Val="  -0.19%  "},
    @{Addr="E29"; Val="  -0.92%  "},
    @{Addr="D30"; Val="'133.37"},
    @{Addr="E30"; Val="  -3.05%  "},
    @{Addr="D31"; Val="1.812.94"},
    @{Addr="E31"; Val="  -1.71%  "},
    @{Addr="D32"; Val="'2.189"},
    @{Addr="E32"; Val="  -5.74%  "},
    @{Addr="D33"; Val="'6.623"},
    @{Addr="E33"; Val="  -6.64%  "},
    @{Addr="D34"; Val="'1.054"},
    @{Addr="E34"; Val="  +6.97%  "},
    @{Addr="D35"; Val="'11.31"},
    @{Addr="E35"; Val="  +1.98%  "},
    @{Addr="D36"; Val="'0.02751"},
    @{Addr="E36"; Val="  -5.54%  "},
    @{Addr="D37"; Val="'0.2493"},
    @{Addr="E37"; Val="  -4.08%  "},
    @{Addr="D38"; Val="'0.08776"},
    @{Addr="E38"; Val="  -1.68%  "},
    @{Addr="D39"; Val="'0.07112"},
    @{Addr="E39"; Val="  -3.71%  "},
    @{Addr="D40"; Val="'6.014"},
    @{Addr="E40"; Val="  -6.56%  "},
    @{Addr="D41"; Val="'0.7063"},
    @{Addr="E41"; Val="  -2.81%  "},
    @{Addr="E42"; Val="  -3.31%  "},
    @{Addr="B43"; Val="EnergySwap"},
    @{Addr="C43"; Val="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"},
    @{Addr="D43"; Val="'15.76"},
    @{Addr="E43"; Val="  -5.75%  "},
    @{Addr="B44"; Val="Aptos"},
    @{Addr="C44"; Val="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"},
    @{Addr="D44"; Val="'12.12"},
    @{Addr="E44"; Val="  -5.35%  "},
    @{Addr="D45"; Val="'0.6531"},
    @{Addr="E45"; Val="  -2.27%  "},
    @{Addr="D46"; Val="'0.9992"},
    @{Addr="E46"; Val="  +0.23%  "},
    @{Addr="D47"; Val="'2.281"},
    @{Addr="E47"; Val="  -5.11%  "},
    @{Addr="D48"; Val="'3.967"},
    @{Addr="E48"; Val="  -1.76%  "},
    @{Addr="D49"; Val="'0.07981"},
    @{Addr="E49"; Val="  -1.42%  "},
    @{Addr="D50"; Val="'127.63"},
    @{Addr="E50"; Val="  -1.66%  "},
    @{Addr="D51"; Val="'1.194"},
    @{Addr="E51"; Val="  -3.77%  "}
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}